# Apply the recorded edits to the workbook:
#  1. Rename the worksheet from "Tabelle1" to "regionescomida"
#  2. Scroll the sheet view back to the top-left and move the selection to E1
#  3. Widen column D slightly (16.83203125 -> 19.33203125)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "regionescomida"

# 2. Reset the scroll position and selection
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E1").Select()

# 3. Widen column D (16.83203125 -> 19.33203125 "characters" in the saved
#    XML). The engine quantizes ColumnWidth to whole pixels, so feed it the
#    COM-units value whose rounded result lands on the recorded width.
$ws.Range("D1").ColumnWidth = 18.5
